$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 (pushes old row 21..152 down to 22..153),
# then fill the new row 21 with a duplicate of the prior row-21 record
# (same market/category/quality/unit/origin/classification) but with its
# own date and volume figure.
$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = 8
$ws.Range("B21").Value = "Terminal La Palmera de La Serena"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44462
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 100112012
$ws.Range("G21").Value = "Espinaca"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 3000
$ws.Range("K21").Value = 400
$ws.Range("L21").Value = 500
$ws.Range("M21").Value = 450
$ws.Range("N21").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O21").Value = "Provincia del Elquí"
$ws.Range("P21").Value = 900
$ws.Range("Q21").Value = 0.5
$ws.Range("R21").Value = "Hortaliza"
